$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.543.11"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "2.603.36"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "3.062.47"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").Value = "59.469.13"
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").Value = "2.626.95"
$ws.Range("E16").Value = "  +1.83%  "
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "341.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("E19").Value = "  +1.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.76%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("E24").Value = "  +1.40%  "
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("D28").Value = "0.0₃0745"
$ws.Range("E28").Value = "  +2.95%  "
$ws.Range("E30").Value = "  +6.15%  "
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.847"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.32%  "
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.825"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "273.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0951"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("E46").Value = "  +3.49%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.940.58"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0223"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("E51").Value = "  +2.06%  "
